# Updates cryptos list (coin prices / 1h volume changes) to match the
# latest scrape. Re-applies each changed cell's new text. Values in the
# "Price" column are stored as literal text (they use '.' as a thousands
# separator in some rows, e.g. "69.105.91"), so a leading apostrophe is
# added whenever the new value would otherwise be auto-parsed by Excel as
# a number (which would silently drop meaningful trailing zeros, e.g.
# "1.00" -> 1). The "Volume(1h)" column values already carry padding
# spaces and a trailing "%" so they are never misread as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    if ($val -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $val
    } else {
        $range.Value = $val
    }
}

Set-TextValue $ws.Range("D2") "69.105.91"
Set-TextValue $ws.Range("E2") "  +2.28%  "
Set-TextValue $ws.Range("D3") "3.817.53"
Set-TextValue $ws.Range("E3") "  +0.98%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "629.80"
Set-TextValue $ws.Range("E5") "  +5.32%  "
Set-TextValue $ws.Range("D6") "165.55"
Set-TextValue $ws.Range("E6") "  +0.57%  "
Set-TextValue $ws.Range("D7") "3.814.23"
Set-TextValue $ws.Range("E7") "  +0.96%  "
Set-TextValue $ws.Range("E8") "  -0.10%  "
Set-TextValue $ws.Range("E9") "  +0.89%  "
Set-TextValue $ws.Range("E10") "  +2.55%  "
Set-TextValue $ws.Range("E11") "  +1.03%  "
Set-TextValue $ws.Range("D12") "6.61"
Set-TextValue $ws.Range("E12") "  +3.13%  "
Set-TextValue $ws.Range("E13") "  +0.85%  "
Set-TextValue $ws.Range("D14") "36.05"
Set-TextValue $ws.Range("E14") "  +1.46%  "
Set-TextValue $ws.Range("D15") "4.456.73"
Set-TextValue $ws.Range("E15") "  +0.93%  "
Set-TextValue $ws.Range("D16") "3.842.17"
Set-TextValue $ws.Range("E16") "  +1.85%  "
Set-TextValue $ws.Range("D17") "69.089.69"
Set-TextValue $ws.Range("E17") "  +2.19%  "
Set-TextValue $ws.Range("D18") "18.04"
Set-TextValue $ws.Range("E18") "  -1.35%  "
Set-TextValue $ws.Range("E19") "  +1.46%  "
Set-TextValue $ws.Range("E20") "  -0.07%  "
Set-TextValue $ws.Range("D21") "465.99"
Set-TextValue $ws.Range("E21") "  +1.22%  "
Set-TextValue $ws.Range("D22") "9.69"
Set-TextValue $ws.Range("E22") "  -0.43%  "
Set-TextValue $ws.Range("D23") "0.710"
Set-TextValue $ws.Range("E23") "  +2.20%  "
Set-TextValue $ws.Range("D24") "0.0000153"
Set-TextValue $ws.Range("E24") "  +4.75%  "
Set-TextValue $ws.Range("D25") "83.67"
Set-TextValue $ws.Range("E25") "  +1.54%  "
Set-TextValue $ws.Range("D26") "12.00"
Set-TextValue $ws.Range("E26") "  +0.17%  "
Set-TextValue $ws.Range("D27") "2.16"
Set-TextValue $ws.Range("E27") "  +3.27%  "
Set-TextValue $ws.Range("D28") "10.05"
Set-TextValue $ws.Range("E28") "  +0.92%  "
Set-TextValue $ws.Range("E29") "  +0.06%  "
Set-TextValue $ws.Range("D30") "3.965.97"
Set-TextValue $ws.Range("E30") "  +0.95%  "
Set-TextValue $ws.Range("E31") "  +1.83%  "
Set-TextValue $ws.Range("E32") "  +1.69%  "
Set-TextValue $ws.Range("E33") "  -2.31%  "
Set-TextValue $ws.Range("D34") "29.21"
Set-TextValue $ws.Range("E34") "  +0.66%  "
Set-TextValue $ws.Range("E35") "  +0.19%  "
Set-TextValue $ws.Range("D36") "9.10"
Set-TextValue $ws.Range("E37") "  +3.97%  "
Set-TextValue $ws.Range("D38") "0.150"
Set-TextValue $ws.Range("E38") "  +8.51%  "
Set-TextValue $ws.Range("D39") "3.44"
Set-TextValue $ws.Range("E39") "  +6.77%  "
Set-TextValue $ws.Range("E40") "  +2.85%  "
Set-TextValue $ws.Range("D41") "0.983"
Set-TextValue $ws.Range("E41") "  -0.33%  "
Set-TextValue $ws.Range("D42") "1.00"
Set-TextValue $ws.Range("E42") "  +0.07%  "
Set-TextValue $ws.Range("E43") "  +0.01%  "
Set-TextValue $ws.Range("D44") "158.02"
Set-TextValue $ws.Range("E44") "  +4.15%  "
Set-TextValue $ws.Range("E45") "  +5.55%  "
Set-TextValue $ws.Range("D46") "0.300"
Set-TextValue $ws.Range("E46") "  +1.26%  "
Set-TextValue $ws.Range("D47") "46.94"
Set-TextValue $ws.Range("E47") "  -1.10%  "
Set-TextValue $ws.Range("B48") "Cosmos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "8.46"
Set-TextValue $ws.Range("E48") "  +1.78%  "
Set-TextValue $ws.Range("B49") "Arweave"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D49") "42.50"
Set-TextValue $ws.Range("E49") "  -1.99%  "
Set-TextValue $ws.Range("B50") "Stacks"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "1.90"
Set-TextValue $ws.Range("E50") "  +2.97%  "
Set-TextValue $ws.Range("D51") "0.000280"
Set-TextValue $ws.Range("E51") "  +13.76%  "
